# specialisatie ipv afgeleid van bij TGs
#
# On the "Standaarden" sheet, the "Inspire-MD-TG-v1.3" (row 8) and
# "Inspire-MD-TG-v2.0" (row 9) records had their relation to
# "ISO-19115:2003" recorded in the "afgeleidvan" column (J). This
# relationship should instead be recorded as a specialization, i.e.
# moved to the "specialisatievan" column (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Standaarden")

# Row 8: Inspire-MD-TG-v1.3
$ws.Range("I8").Value2 = $ws.Range("J8").Value2
$ws.Range("J8").ClearContents()

# Row 9: Inspire-MD-TG-v2.0
$ws.Range("I9").Value2 = $ws.Range("J9").Value2
$ws.Range("J9").ClearContents()
